$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.99"
$ws.Range("E2").Value = "'8.79%"
$ws.Range("D3").Value = "'49.53"
$ws.Range("E3").Value = "'18.39%"
$ws.Range("D4").Value = "'5.350"
$ws.Range("E4").Value = "'6.93%"
$ws.Range("D5").Value = "'0.08153"
$ws.Range("E5").Value = "'8.47%"
$ws.Range("D6").Value = "'4.606"
$ws.Range("E6").Value = "'5.15%"
$ws.Range("D7").Value = "'1.677"
$ws.Range("E7").Value = "'5.90%"
$ws.Range("D8").Value = "'1.163"
$ws.Range("E8").Value = "'25.54%"
$ws.Range("D9").Value = "'0.1352"
$ws.Range("E9").Value = "'12.86%"
$ws.Range("D10").Value = "'0.1960"
$ws.Range("E10").Value = "'7.16%"
$ws.Range("D11").Value = "'0.09551"
$ws.Range("E11").Value = "'7.67%"
$ws.Range("D12").Value = "'0.04570"
$ws.Range("E12").Value = "'12.19%"
$ws.Range("E13").Value = "'-0.07%"
$ws.Range("D14").Value = "'0.001338"
$ws.Range("E14").Value = "'4.16%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04298"
$ws.Range("E15").Value = "'4.80%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005855"
$ws.Range("E16").Value = "'-1.49%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.385"
$ws.Range("E17").Value = "'0.85%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.440"
$ws.Range("E18").Value = "'1.62%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3394"
$ws.Range("E19").Value = "'2.45%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'8.117"
$ws.Range("E20").Value = "'-0.15%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1409"
$ws.Range("E21").Value = "'1.36%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.3052"
$ws.Range("E22").Value = "'-1.60%"
$ws.Range("E23").Value = "'3.25%"
$ws.Range("D24").Value = "'0.004309"
$ws.Range("E24").Value = "'10.30%"
$ws.Range("E25").Value = "'9.65%"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("D38").Value = "'0.02770"
$ws.Range("E38").Value = "'15.19%"
$ws.Range("D39").Value = "'0.05531"
$ws.Range("E39").Value = "'5.46%"
$ws.Range("D41").Value = "'0.007751"
$ws.Range("E41").Value = "'-0.64%"
$ws.Range("E42").Value = "'9.14%"
$ws.Range("D43").Value = "'0.007685"
$ws.Range("E43").Value = "'4.02%"
$ws.Range("D44").Value = "'0.008098"
$ws.Range("E44").Value = "'11.89%"
$ws.Range("D45").Value = "'0.3508"
$ws.Range("E45").Value = "'18.63%"
$ws.Range("D46").Value = "'0.00006769"
$ws.Range("E46").Value = "'3.71%"
$ws.Range("E48").Value = "'18.04%"
$ws.Range("E49").Value = "'-4.83%"